$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-formatted cells keep their literal string representation
# (Price/Volume columns contain numeric-looking strings like "249.07" or
# "2.00" that Excel would otherwise auto-convert to numbers on assignment).

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '41.012.24'
$ws.Range("E2").NumberFormat = "@"
$ws.Cells.Item(2, 5).Value = '  -2.09%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '2.171.98'
$ws.Range("E3").NumberFormat = "@"
$ws.Cells.Item(3, 5).Value = '  -3.07%  '

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Cells.Item(4, 5).Value = '  -0.05%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '249.07'
$ws.Range("E5").NumberFormat = "@"
$ws.Cells.Item(5, 5).Value = '  -1.32%  '

# Row 6
$ws.Range("E6").NumberFormat = "@"
$ws.Cells.Item(6, 5).Value = '  -2.45%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '66.93'
$ws.Range("E7").NumberFormat = "@"
$ws.Cells.Item(7, 5).Value = '  -7.13%  '

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Cells.Item(8, 5).Value = '  +0.10%  '

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Cells.Item(9, 5).Value = '  -0.20%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '58.46'
$ws.Range("E10").NumberFormat = "@"
$ws.Cells.Item(10, 5).Value = '  -0.05%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '35.74'
$ws.Range("E11").NumberFormat = "@"
$ws.Cells.Item(11, 5).Value = '  -15.55%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.0925'
$ws.Range("E12").NumberFormat = "@"
$ws.Cells.Item(12, 5).Value = '  -5.24%  '

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Cells.Item(13, 5).Value = '  -1.66%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '6.88'
$ws.Range("E14").NumberFormat = "@"
$ws.Cells.Item(14, 5).Value = '  -0.56%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '2.495.44'
$ws.Range("E15").NumberFormat = "@"
$ws.Cells.Item(15, 5).Value = '  -3.01%  '

# Row 16
$ws.Range("E16").NumberFormat = "@"
$ws.Cells.Item(16, 5).Value = '  -0.01%  '

# Row 17
$ws.Range("E17").NumberFormat = "@"
$ws.Cells.Item(17, 5).Value = '  -6.62%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '2.176.00'
$ws.Range("E18").NumberFormat = "@"
$ws.Cells.Item(18, 5).Value = '  -2.92%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '40.869.89'
$ws.Range("E19").NumberFormat = "@"
$ws.Cells.Item(19, 5).Value = '  -2.28%  '

# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Cells.Item(20, 5).Value = '  -3.42%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '6.08'

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '71.37'
$ws.Range("E22").NumberFormat = "@"
$ws.Cells.Item(22, 5).Value = '  -2.93%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '229.64'
$ws.Range("E23").NumberFormat = "@"
$ws.Cells.Item(23, 5).Value = '  -2.59%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '2.07'
$ws.Range("E24").NumberFormat = "@"
$ws.Cells.Item(24, 5).Value = '  -8.15%  '

# Row 26
$ws.Cells.Item(26, 2).Value = 'WEMIXToken'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D26").NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '3.74'
$ws.Range("E26").NumberFormat = "@"
$ws.Cells.Item(26, 5).Value = '  -0.62%  '

# Row 27
$ws.Cells.Item(27, 2).Value = 'Cosmos'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '11.35'
$ws.Range("E27").NumberFormat = "@"
$ws.Cells.Item(27, 5).Value = '  +11.33%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '2.42'
$ws.Range("E28").NumberFormat = "@"
$ws.Cells.Item(28, 5).Value = '  -3.55%  '

# Row 29
$ws.Cells.Item(29, 2).Value = 'Monero'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D29").NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '167.78'
$ws.Range("E29").NumberFormat = "@"
$ws.Cells.Item(29, 5).Value = '  -2.62%  '

# Row 30
$ws.Cells.Item(30, 2).Value = 'Toncoin'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '2.00'
$ws.Range("E30").NumberFormat = "@"
$ws.Cells.Item(30, 5).Value = '  -11.75%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '20.15'
$ws.Range("E31").NumberFormat = "@"
$ws.Cells.Item(31, 5).Value = '  -3.13%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '0.122'
$ws.Range("E32").NumberFormat = "@"
$ws.Cells.Item(32, 5).Value = '  -1.03%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '5.63'
$ws.Range("E33").NumberFormat = "@"
$ws.Cells.Item(33, 5).Value = '  +2.83%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.0739'
$ws.Range("E34").NumberFormat = "@"
$ws.Cells.Item(34, 5).Value = '  +2.09%  '

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Cells.Item(35, 5).Value = '  -3.59%  '

# Row 36
$ws.Range("E36").NumberFormat = "@"
$ws.Cells.Item(36, 5).Value = '  -3.84%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '4.07'
$ws.Range("E37").NumberFormat = "@"
$ws.Cells.Item(37, 5).Value = '  -2.21%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '25.49'
$ws.Range("E38").NumberFormat = "@"
$ws.Cells.Item(38, 5).Value = '  -4.15%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.0301'
$ws.Range("E39").NumberFormat = "@"
$ws.Cells.Item(39, 5).Value = '  +7.53%  '

# Row 40
$ws.Range("E40").NumberFormat = "@"
$ws.Cells.Item(40, 5).Value = '  -5.73%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '5.50'
$ws.Range("E41").NumberFormat = "@"
$ws.Cells.Item(41, 5).Value = '  -9.89%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '11.57'
$ws.Range("E42").NumberFormat = "@"
$ws.Cells.Item(42, 5).Value = '  -0.54%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '61.26'
$ws.Range("E43").NumberFormat = "@"
$ws.Cells.Item(43, 5).Value = '  -11.64%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '4.78'
$ws.Range("E44").NumberFormat = "@"
$ws.Cells.Item(44, 5).Value = '  -5.68%  '

# Row 45
$ws.Range("E45").NumberFormat = "@"
$ws.Cells.Item(45, 5).Value = '  -11.22%  '

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Cells.Item(46, 5).Value = '  -4.56%  '

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Cells.Item(47, 5).Value = '  +0.50%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '1.17'
$ws.Range("E48").NumberFormat = "@"
$ws.Cells.Item(48, 5).Value = '  +3.36%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.0982'
$ws.Range("E49").NumberFormat = "@"
$ws.Cells.Item(49, 5).Value = '  -3.96%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '1.15'
$ws.Range("E50").NumberFormat = "@"
$ws.Cells.Item(50, 5).Value = '  -4.06%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '2.71'
$ws.Range("E51").NumberFormat = "@"
$ws.Cells.Item(51, 5).Value = '  -0.56%  '
